$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.652.42"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "2.069.87"
$ws.Range("E3").Value = "  -1.73%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.75"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.77"
$ws.Range("E8").Value = "  -0.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("D12").Value = "2.373.31"
$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.74"
$ws.Range("E13").Value = "  +1.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.74"
$ws.Range("E14").Value = "  -2.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.769"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").Value = "2.093.65"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "37.582.77"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("E19").Value = "  -3.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.89"
$ws.Range("E20").Value = "  +1.07%  "

$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.02"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  -0.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.20"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.138"
$ws.Range("E27").Value = "  +3.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.96"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.34"
$ws.Range("E29").Value = "  -1.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  -2.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +1.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.65"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0627"
$ws.Range("E33").Value = "  +0.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.62"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  -5.46%  "

$ws.Range("E36").Value = "  +2.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.37"
$ws.Range("E37").Value = "  -3.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.29"
$ws.Range("E39").Value = "  -4.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0967"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.62"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("E42").Value = "  -2.23%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "1.454.61"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.36"
$ws.Range("E45").Value = "  +5.22%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.15"
$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.17"
$ws.Range("E47").Value = "  -1.68%  "

$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.01"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("D51").Value = "2.257.35"
$ws.Range("E51").Value = "  -2.02%  "
